$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "cost" (C1) to "euro"
$ws.Range("C1").Value = "euro"

# Update the selection to be on cell C1 (single cell selection)
$ws.Range("C1").Select() | Out-Null
